$d = $word.ActiveDocument

# --- "Day" -> "Time" -----------------------------------------------------
# The two "Day: day of year in which a census was taken" bullets get their
# leading word changed from "Day" to "Time". When this is done narrowly in
# Word (selecting just "Day" and typing over it) the run that used to hold
# "Day: day of year..." is split into two runs at the edit boundary: one
# holding "Time" and one holding ": day of year in which a census was
# taken" - both keep identical run formatting.
while ($true) {
    $rng = $d.Content
    $found = $rng.Find.Execute("Day", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }

    # Touch ascii/hAnsi font first...
    $fontName = $rng.Font.Name
    $rng.Font.Name = $fontName
    $rng.Text = "Time"

    # ...then re-touch the complex-script font on the just-replaced word
    # only (Find.Execute re-seats $rng to span exactly "Time" after the
    # assignment above). Re-seating this run property keeps the freshly
    # split run distinct rather than letting it silently re-merge with its
    # sibling.
    $rng.Font.NameBi = $fontName
}

# --- "Adults" -> "A" ------------------------------------------------------
# The "Adults: adult density..." bullets simply get their leading word
# shortened to "A"; no run split occurs here. Match on "Adults:" (with the
# trailing colon) so the unrelated "(Adults - A_SE)" / "(Adults + A_SE)"
# mentions further down are left untouched.
while ($d.Content.Find.Execute("Adults:", $true, $false, $false, $false, $false, $true, 1, $false, "A:", 2)) {
}
